# Fruta / hortaliza, semanal
# Insert two new weekly price records (rows 79-80) into the "Poroto granado"
# sheet, pushing the existing rows 79-107 down to rows 81-109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 79.
$ws.Rows("79:80").Insert()

# New row 79: Region de La Araucania, 07-02-2023, Primera, 50 sacos @ 45000
$ws.Range("A79").Value = 10
$ws.Range("B79").Value = "Vega Modelo de Temuco"
$ws.Range("C79").Value = "La Araucanía"
$ws.Range("D79").Value = 44964
$ws.Range("E79").Value = 9
$ws.Range("F79").Value = 100112030
$ws.Range("G79").Value = "Poroto granado"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 50
$ws.Range("K79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("M79").Value = 45000
$ws.Range("N79").Value = "`$/saco 25 kilos"
$ws.Range("O79").Value = "Región de La Araucanía"
$ws.Range("P79").Value = 1800
$ws.Range("Q79").Value = 25
$ws.Range("R79").Value = "Hortaliza"

# New row 80: Region del Maule, 07-02-2023, Primera, 30 sacos @ 45000
$ws.Range("A80").Value = 10
$ws.Range("B80").Value = "Vega Modelo de Temuco"
$ws.Range("C80").Value = "La Araucanía"
$ws.Range("D80").Value = 44964
$ws.Range("E80").Value = 9
$ws.Range("F80").Value = 100112030
$ws.Range("G80").Value = "Poroto granado"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 30
$ws.Range("K80").Value = 45000
$ws.Range("L80").Value = 45000
$ws.Range("M80").Value = 45000
$ws.Range("N80").Value = "`$/saco 25 kilos"
$ws.Range("O80").Value = "Región del Maule"
$ws.Range("P80").Value = 1800
$ws.Range("Q80").Value = 25
$ws.Range("R80").Value = "Hortaliza"
